$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.859.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.52%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.915.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.96%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.86%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.506'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.914.98'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.99%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.80%  '

$ws.Range('E11').Value = '  +4.63%  '

$ws.Range('E12').Value = '  -2.97%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000239'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.10%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.63%  '

$ws.Range('E15').Value = '  -1.56%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.395.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.00%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.809.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.89%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.910.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.23%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '436.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.22%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.62%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.660'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.47%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.10%  '

$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.53%  '

$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000108'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +22.24%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.30%  '

$ws.Range('E31').Value = '  -2.11%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.80%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.110'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.84%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '

$ws.Range('E35').Value = '  -0.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.977'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.57%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.34%  '

$ws.Range('E39').Value = '  -1.05%  '

$ws.Range('E40').Value = '  -0.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.273'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.697.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.23%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.29%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0338'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '348.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.29%  '

$ws.Range('E49').Value = '  +0.02%  '

$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.88%  '
